# Generate Report for Handoff
# Adds two new source files (64b700f9... and 6f40ab55...) to the
# localization-status workbook, ahead of the ".localization-config" row,
# across all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": columns A=File Name, B=zh-cn, C=de-de
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Duplicate the ".localization-config" row (row 4) twice so the two new
# rows inherit its formatting (hyperlink-styled column A, etc.), pushing
# the original row down to row 6.
$ws.Rows.Item(4).Copy()
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Copy()
$ws.Rows.Item(4).Insert()

# Row 4: 64b700f9...md
$ws.Range("A4").Value = "64b700f9-f2ed-4a84-84d9-100e4bd1efaa.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"

# Row 5: 6f40ab55...md
$ws.Range("A5").Value = "6f40ab55-d572-4e31-b14c-0f4587129f02.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"

# Row 6 already holds the original ".localization-config" / "Not to be
# localized" values (shifted down automatically) - nothing else to do.

# Rebuild hyperlinks in final order (this engine does not auto-shift
# hyperlink refs on row insert).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1e28b551f0e8fd33d9b8307f67d8bad6ea52e31d/e2e/aaacb84e-c7ef-45a9-94bb-7e392c7f5142.md", "", "", "aaacb84e-c7ef-45a9-94bb-7e392c7f5142.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/1e28b551f0e8fd33d9b8307f67d8bad6ea52e31d/e2e/bf3a1c35-206d-42ae-b422-1131219ec14d.md", "", "", "bf3a1c35-206d-42ae-b422-1131219ec14d.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/1e28b551f0e8fd33d9b8307f67d8bad6ea52e31d/e2e/64b700f9-f2ed-4a84-84d9-100e4bd1efaa.md", "", "", "64b700f9-f2ed-4a84-84d9-100e4bd1efaa.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/1e28b551f0e8fd33d9b8307f67d8bad6ea52e31d/e2e/6f40ab55-d572-4e31-b14c-0f4587129f02.md", "", "", "6f40ab55-d572-4e31-b14c-0f4587129f02.md")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/1e28b551f0e8fd33d9b8307f67d8bad6ea52e31d/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "zh-cn": A=Source File Name, B=Status, C=Latest Handoff File,
# D=Latest Handoff Datetime, E=Latest Target File, F=Latest Handback
# File, G=Latest Handback DateTime, H=Handoff Reason, I=Dependency From
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 3 already has every column we need (A,B,C,D,G,H) with the right
# styles, so duplicate it twice ahead of the ".localization-config" row
# (row 4), pushing the original down to row 6.
$ws.Rows.Item(3).Copy()
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(3).Copy()
$ws.Rows.Item(4).Insert()

# Row 4: 64b700f9...md handoff
$ws.Range("A4").Value = "64b700f9-f2ed-4a84-84d9-100e4bd1efaa.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "64b700f9-f2ed-4a84-84d9-100e4bd1efaa.c47bbf2313f119567aea32e68cb2cbcc107fe1af.zh-cn.xlf"
$ws.Range("D4").Value = "2016-02-18 07:28:29"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Include"

# Row 5: 6f40ab55...md handoff
$ws.Range("A5").Value = "6f40ab55-d572-4e31-b14c-0f4587129f02.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "6f40ab55-d572-4e31-b14c-0f4587129f02.54be0bc29226a42ecb890ea7fbd30325634e2456.zh-cn.xlf"
$ws.Range("D5").Value = "2016-02-18 07:28:29"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Include"

# Row 6 already holds the original ".localization-config" values.

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1e28b551f0e8fd33d9b8307f67d8bad6ea52e31d/e2e/aaacb84e-c7ef-45a9-94bb-7e392c7f5142.md", "", "", "aaacb84e-c7ef-45a9-94bb-7e392c7f5142.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c85ab8eb8e51f71e076a9167273184c92ae5ddf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/aaacb84e-c7ef-45a9-94bb-7e392c7f5142.e617407039e4321f931ebeedbff092c377d24cd3.zh-cn.xlf", "", "", "aaacb84e-c7ef-45a9-94bb-7e392c7f5142.e617407039e4321f931ebeedbff092c377d24cd3.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/1e28b551f0e8fd33d9b8307f67d8bad6ea52e31d/e2e/bf3a1c35-206d-42ae-b422-1131219ec14d.md", "", "", "bf3a1c35-206d-42ae-b422-1131219ec14d.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c85ab8eb8e51f71e076a9167273184c92ae5ddf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/bf3a1c35-206d-42ae-b422-1131219ec14d.daa023d876985a6e4d82a761633242467fd31e69.zh-cn.xlf", "", "", "bf3a1c35-206d-42ae-b422-1131219ec14d.daa023d876985a6e4d82a761633242467fd31e69.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/1e28b551f0e8fd33d9b8307f67d8bad6ea52e31d/e2e/64b700f9-f2ed-4a84-84d9-100e4bd1efaa.md", "", "", "64b700f9-f2ed-4a84-84d9-100e4bd1efaa.md")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c85ab8eb8e51f71e076a9167273184c92ae5ddf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/64b700f9-f2ed-4a84-84d9-100e4bd1efaa.c47bbf2313f119567aea32e68cb2cbcc107fe1af.zh-cn.xlf", "", "", "64b700f9-f2ed-4a84-84d9-100e4bd1efaa.c47bbf2313f119567aea32e68cb2cbcc107fe1af.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/1e28b551f0e8fd33d9b8307f67d8bad6ea52e31d/e2e/6f40ab55-d572-4e31-b14c-0f4587129f02.md", "", "", "6f40ab55-d572-4e31-b14c-0f4587129f02.md")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c85ab8eb8e51f71e076a9167273184c92ae5ddf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/6f40ab55-d572-4e31-b14c-0f4587129f02.54be0bc29226a42ecb890ea7fbd30325634e2456.zh-cn.xlf", "", "", "6f40ab55-d572-4e31-b14c-0f4587129f02.54be0bc29226a42ecb890ea7fbd30325634e2456.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/1e28b551f0e8fd33d9b8307f67d8bad6ea52e31d/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "de-de": same layout as "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(3).Copy()
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(3).Copy()
$ws.Rows.Item(4).Insert()

# Row 4: 64b700f9...md handoff
$ws.Range("A4").Value = "64b700f9-f2ed-4a84-84d9-100e4bd1efaa.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "64b700f9-f2ed-4a84-84d9-100e4bd1efaa.c47bbf2313f119567aea32e68cb2cbcc107fe1af.de-de.xlf"
$ws.Range("D4").Value = "2016-02-18 07:28:40"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Include"

# Row 5: 6f40ab55...md handoff
$ws.Range("A5").Value = "6f40ab55-d572-4e31-b14c-0f4587129f02.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "6f40ab55-d572-4e31-b14c-0f4587129f02.54be0bc29226a42ecb890ea7fbd30325634e2456.de-de.xlf"
$ws.Range("D5").Value = "2016-02-18 07:28:40"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Include"

# Row 6 already holds the original ".localization-config" values.

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1e28b551f0e8fd33d9b8307f67d8bad6ea52e31d/e2e/aaacb84e-c7ef-45a9-94bb-7e392c7f5142.md", "", "", "aaacb84e-c7ef-45a9-94bb-7e392c7f5142.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/014eb1c6d68c65332d16fc7721e1b77b274d8045/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/aaacb84e-c7ef-45a9-94bb-7e392c7f5142.e617407039e4321f931ebeedbff092c377d24cd3.de-de.xlf", "", "", "aaacb84e-c7ef-45a9-94bb-7e392c7f5142.e617407039e4321f931ebeedbff092c377d24cd3.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/1e28b551f0e8fd33d9b8307f67d8bad6ea52e31d/e2e/bf3a1c35-206d-42ae-b422-1131219ec14d.md", "", "", "bf3a1c35-206d-42ae-b422-1131219ec14d.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/014eb1c6d68c65332d16fc7721e1b77b274d8045/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/bf3a1c35-206d-42ae-b422-1131219ec14d.daa023d876985a6e4d82a761633242467fd31e69.de-de.xlf", "", "", "bf3a1c35-206d-42ae-b422-1131219ec14d.daa023d876985a6e4d82a761633242467fd31e69.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/1e28b551f0e8fd33d9b8307f67d8bad6ea52e31d/e2e/64b700f9-f2ed-4a84-84d9-100e4bd1efaa.md", "", "", "64b700f9-f2ed-4a84-84d9-100e4bd1efaa.md")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/014eb1c6d68c65332d16fc7721e1b77b274d8045/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/64b700f9-f2ed-4a84-84d9-100e4bd1efaa.c47bbf2313f119567aea32e68cb2cbcc107fe1af.de-de.xlf", "", "", "64b700f9-f2ed-4a84-84d9-100e4bd1efaa.c47bbf2313f119567aea32e68cb2cbcc107fe1af.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/1e28b551f0e8fd33d9b8307f67d8bad6ea52e31d/e2e/6f40ab55-d572-4e31-b14c-0f4587129f02.md", "", "", "6f40ab55-d572-4e31-b14c-0f4587129f02.md")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/014eb1c6d68c65332d16fc7721e1b77b274d8045/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/6f40ab55-d572-4e31-b14c-0f4587129f02.54be0bc29226a42ecb890ea7fbd30325634e2456.de-de.xlf", "", "", "6f40ab55-d572-4e31-b14c-0f4587129f02.54be0bc29226a42ecb890ea7fbd30325634e2456.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/1e28b551f0e8fd33d9b8307f67d8bad6ea52e31d/.localization-config", "", "", ".localization-config")

Write-Output "Report generated for handoff."
